$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Range("D5").Value = "Obesity"
